$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.742267333333333
$ws.Range("H2").Value = 5.226802
$ws.Range("I2").Value = 0.02937931807530601
$ws.Range("J2").Value = 0.029379318075306
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.171693666666666
$ws.Range("N2").Value = 12.515081
$ws.Range("O2").Value = 0.1077921033402881
$ws.Range("P2").Value = 0.1077921033402881
$ws.Range("Q2").Value = 7.268205600106889
$ws.Range("R2").Value = 65.41385040096199
$ws.Range("S2").Value = 0.00316685849004058
$ws.Range("T2").Value = 0.00316685849004058
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.742267333333333
$ws.Range("H3").Value = 5.226802
$ws.Range("I3").Value = 0.02937931807530601
$ws.Range("J3").Value = 0.029379318075306
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 20.39394566666667
$ws.Range("N3").Value = 61.181837
$ws.Range("O3").Value = 0.5269577477327285
$ws.Range("P3").Value = 0.5269577477327286
$ws.Range("Q3").Value = 35.53170533280822
$ws.Range("R3").Value = 319.785347995274
$ws.Range("S3").Value = 0.0154816592828867
$ws.Range("T3").Value = 0.0154816592828867
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.742267333333333
$ws.Range("H4").Value = 5.226802
$ws.Range("I4").Value = 0.02937931807530601
$ws.Range("J4").Value = 0.029379318075306
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.13565266666667
$ws.Range("N4").Value = 42.406958
$ws.Range("O4").Value = 0.3652501489269833
$ws.Range("P4").Value = 0.3652501489269833
$ws.Range("Q4").Value = 24.62808587647956
$ws.Range("R4").Value = 221.652772888316
$ws.Range("S4").Value = 0.01073080030237873
$ws.Range("T4").Value = 0.01073080030237873
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 37.52183533333334
$ws.Range("H5").Value = 112.565506
$ws.Range("I5").Value = 0.632719166534674
$ws.Range("J5").Value = 0.6327191665346739
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.171693666666666
$ws.Range("N5").Value = 12.515081
$ws.Range("O5").Value = 0.1077921033402881
$ws.Range("P5").Value = 0.1077921033402881
$ws.Range("Q5").Value = 156.5296028217762
$ws.Range("R5").Value = 1408.766425395986
$ws.Range("S5").Value = 0.06820212978448655
$ws.Range("T5").Value = 0.06820212978448655
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 37.52183533333334
$ws.Range("H6").Value = 112.565506
$ws.Range("I6").Value = 0.632719166534674
$ws.Range("J6").Value = 0.6327191665346739
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 20.39394566666667
$ws.Range("N6").Value = 61.181837
$ws.Range("O6").Value = 0.5269577477327285
$ws.Range("P6").Value = 0.5269577477327286
$ws.Range("Q6").Value = 765.2182711016136
$ws.Range("R6").Value = 6886.964439914523
$ws.Range("S6").Value = 0.333416266944441
$ws.Range("T6").Value = 0.333416266944441
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 37.52183533333334
$ws.Range("H7").Value = 112.565506
$ws.Range("I7").Value = 0.632719166534674
$ws.Range("J7").Value = 0.6327191665346739
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.13565266666667
$ws.Range("N7").Value = 42.406958
$ws.Range("O7").Value = 0.3652501489269833
$ws.Range("P7").Value = 0.3652501489269833
$ws.Range("Q7").Value = 530.3956316878609
$ws.Range("R7").Value = 4773.560685190749
$ws.Range("S7").Value = 0.2311007698057464
$ws.Range("T7").Value = 0.2311007698057464
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 20.03840833333333
$ws.Range("H8").Value = 60.115225
$ws.Range("I8").Value = 0.3379015153900201
$ws.Range("J8").Value = 0.33790151539002
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.171693666666666
$ws.Range("N8").Value = 12.515081
$ws.Range("O8").Value = 0.1077921033402881
$ws.Range("P8").Value = 0.1077921033402881
$ws.Range("Q8").Value = 83.59410113424721
$ws.Range("R8").Value = 752.3469102082249
$ws.Range("S8").Value = 0.036423115065761
$ws.Range("T8").Value = 0.036423115065761
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 20.03840833333333
$ws.Range("H9").Value = 60.115225
$ws.Range("I9").Value = 0.3379015153900201
$ws.Range("J9").Value = 0.33790151539002
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 20.39394566666667
$ws.Range("N9").Value = 61.181837
$ws.Range("O9").Value = 0.5269577477327285
$ws.Range("P9").Value = 0.5269577477327286
$ws.Range("Q9").Value = 408.6622107964806
$ws.Range("R9").Value = 3677.959897168325
$ws.Range("S9").Value = 0.1780598215054009
$ws.Range("T9").Value = 0.1780598215054009
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 20.03840833333333
$ws.Range("H10").Value = 60.115225
$ws.Range("I10").Value = 0.3379015153900201
$ws.Range("J10").Value = 0.33790151539002
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.13565266666667
$ws.Range("N10").Value = 42.406958
$ws.Range("O10").Value = 0.3652501489269833
$ws.Range("P10").Value = 0.3652501489269833
$ws.Range("Q10").Value = 283.2559801928389
$ws.Range("R10").Value = 2549.30382173555
$ws.Range("S10").Value = 0.1234185788188581
$ws.Range("T10").Value = 0.1234185788188581